$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D16").Value = "CAMERAS: Enhanced Resolution And Sanity Preserving Class Activation Mapping For Image Saliency 내용 정리 [XAI-24]"
$ws.Range("E16").Value = "https://wewinserv.tistory.com/187"

$ws.Range("D32").Value = "파이썬 패키지 개념 (feat. 코딩도장)"
$ws.Range("E32").Value = "https://dodonam.tistory.com/340"

$ws.Range("D41").Value = "신속하고 편리한 개발, Lowcode"
$ws.Range("E41").Value = "http://cloudinsight.net/cloud/%ec%8b%a0%ec%86%8d%ed%95%98%ea%b3%a0-%ed%8e%b8%eb%a6%ac%ed%95%9c-%ea%b0%9c%eb%b0%9c-lowcode/"

$wb.Save()
